$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:E -> B:F
$ws.Columns.Item(1).EntireColumn.Insert()

# Copy the header formatting (bold/border/alignment) from B1 into the new A1
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)

# Set the new "ID" header and identifiers for each row
$ws.Cells.Item(1, 1).Value = "ID"

$ws.Cells.Item(2, 1).Value = "Hb 2"
$ws.Cells.Item(3, 1).Value = "Hb 3"
$ws.Cells.Item(4, 1).Value = "S 24"
$ws.Cells.Item(5, 1).Value = "S 28"
$ws.Cells.Item(6, 1).Value = "Hb 107"
$ws.Cells.Item(7, 1).Value = "Hb 66"
$ws.Cells.Item(8, 1).Value = "Hb 69"
$ws.Cells.Item(9, 1).Value = "Hb 95"
$ws.Cells.Item(10, 1).Value = "Hb 99"
$ws.Cells.Item(11, 1).Value = "Hb 92"
$ws.Cells.Item(12, 1).Value = "Hb 40"
$ws.Cells.Item(13, 1).Value = "Hb 41"
$ws.Cells.Item(14, 1).Value = "S 11"
$ws.Cells.Item(15, 1).Value = "Hb 57"
$ws.Cells.Item(16, 1).Value = "S 21"
$ws.Cells.Item(17, 1).Value = "S 22"
$ws.Cells.Item(18, 1).Value = "S 3"
$ws.Cells.Item(19, 1).Value = "S 4"
$ws.Cells.Item(20, 1).Value = "S 5"
$ws.Cells.Item(21, 1).Value = "Hb 74"
$ws.Cells.Item(22, 1).Value = "Hb 79"
$ws.Cells.Item(23, 1).Value = "Hb 32"
$ws.Cells.Item(24, 1).Value = "S 15"
$ws.Cells.Item(25, 1).Value = "S 16"

$excel.CutCopyMode = 0
